$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test-Cases")

# Update rows 6, 7, 8: Approved/Rejected column (I) becomes "Approved"
# and ReasonToReject column (J) is cleared.
$ws.Range("I6").Value = "Approved"
$ws.Range("J6").ClearContents()

$ws.Range("I7").Value = "Approved"
$ws.Range("J7").ClearContents()

$ws.Range("I8").Value = "Approved"
$ws.Range("J8").ClearContents()

# Update the active selection on the sheet
$ws.Range("H21").Select()
